$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: add a time value (22:00) in column H ---
$ws.Range("H17").Value = 0.91666666666666663
$ws.Range("H17").NumberFormat = "h:mm"

# --- Row 18: new date row ---
$ws.Range("C18").Value = 41971
$ws.Range("C18").NumberFormat = $ws.Range("C17").NumberFormat
$ws.Range("F18").Value = 1

# --- Row 19: new date row ---
$ws.Range("C19").Value = 41972
$ws.Range("C19").NumberFormat = $ws.Range("C17").NumberFormat
$ws.Range("D19").Value = 1
$ws.Range("F19").Value = 2
$ws.Range("H19").NumberFormat = "h:mm"

# --- Row 20: totals ---
$ws.Range("D20").Formula = "=SUM(D17:D19)"
$ws.Range("E20").Formula = "=SUM(E17:E19)"
$ws.Range("F20").Formula = "=SUM(F17:F19)"
$ws.Range("H20").Value = "Total"
$ws.Range("I20").Formula = "=SUM(D20:F20)"

# --- Update the active selection to mirror the authored workbook ---
$ws.Range("F19").Select() | Out-Null
